# Applies the cryptos-list price refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-column (D) values are digit/dot strings that Excel would otherwise
# auto-convert to numbers (dropping trailing zeros, etc). Force them to stay
# text the same way typing an apostrophe-prefixed value does, then drop the
# quote-prefix formatting artifact so the cell style matches the original
# (no explicit style index).
function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

function Set-PlainValue($ref, $val) {
    $ws.Range($ref).Value = $val
}

$priceUpdates = @(
    @("D2", "25.796.55"),
    @("D3", "1.637.56"),
    @("D4", "1.003"),
    @("D5", "215.19"),
    @("D6", "0.5061"),
    @("D8", "0.2578"),
    @("D9", "0.06411"),
    @("D10", "20.36"),
    @("D11", "0.07794"),
    @("D12", "4.249"),
    @("D13", "1.865.60"),
    @("D14", "1.631.72"),
    @("D15", "0.5605"),
    @("D16", "0.0₅7630"),
    @("D17", "63.29"),
    @("D18", "25.813.38"),
    @("D20", "4.375"),
    @("D21", "191.62"),
    @("D22", "9.899"),
    @("D23", "6.140"),
    @("D25", "1.792"),
    @("D26", "139.88"),
    @("D27", "0.1228"),
    @("D28", "6.811"),
    @("D29", "15.52"),
    @("D30", "1.242"),
    @("D31", "0.04939"),
    @("D32", "3.279"),
    @("D34", "1.568"),
    @("D35", "2.385"),
    @("D36", "0.9018"),
    @("D37", "2.571"),
    @("D38", "0.5556"),
    @("D39", "1.132.41"),
    @("D40", "0.01568"),
    @("D41", "1.002"),
    @("D42", "0.8000"),
    @("D43", "5.448"),
    @("D44", "98.80"),
    @("D45", "1.776.81"),
    @("D46", "0.0₈113"),
    @("D47", "55.54"),
    @("D48", "0.4259"),
    @("D49", "7.796"),
    @("D50", "0.05031"),
    @("D51", "0.9983")
)

$otherUpdates = @(
    @("E2", "  -0.46%  "),
    @("E3", "  +0.27%  "),
    @("E4", "  +0.19%  "),
    @("E5", "  -0.40%  "),
    @("E6", "  -1.38%  "),
    @("E8", "  +0.41%  "),
    @("E9", "  +0.96%  "),
    @("E10", "  +4.44%  "),
    @("E11", "  +0.24%  "),
    @("E12", "  -0.13%  "),
    @("B13", "WrappedliquidstakedEther2.0"),
    @("C13", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"),
    @("E13", "  +0.38%  "),
    @("B14", "WrappedEther"),
    @("C14", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"),
    @("E14", "  -0.15%  "),
    @("E15", "  +1.48%  "),
    @("E16", "  -0.01%  "),
    @("E17", "  -0.84%  "),
    @("E18", "  -0.48%  "),
    @("E19", "  +0.21%  "),
    @("E20", "  -1.17%  "),
    @("E21", "  -2.01%  "),
    @("E22", "  +0.35%  "),
    @("E23", "  +1.84%  "),
    @("E24", "  +0.29%  "),
    @("E25", "  -5.38%  "),
    @("E26", "  -1.73%  "),
    @("E27", "  -2.79%  "),
    @("E28", "  +0.69%  "),
    @("E29", "  -0.68%  "),
    @("E30", "  +0.00%  "),
    @("E31", "  +0.44%  "),
    @("E32", "  +1.44%  "),
    @("E33", "  +1.87%  "),
    @("E34", "  +1.45%  "),
    @("E35", "  +0.50%  "),
    @("E36", "  +0.38%  "),
    @("E37", "  +1.38%  "),
    @("E38", "  +0.58%  "),
    @("E39", "  +1.51%  "),
    @("E40", "  +0.56%  "),
    @("E41", "  +0.17%  "),
    @("B42", "TrustWalletToken"),
    @("C42", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"),
    @("E42", "  +0.22%  "),
    @("B43", "FraxShare"),
    @("C43", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"),
    @("E43", "  -2.39%  "),
    @("E44", "  +1.05%  "),
    @("E45", "  +0.49%  "),
    @("E46", "  -3.65%  "),
    @("E47", "  +1.32%  "),
    @("E48", "  -3.96%  "),
    @("E49", "  +3.14%  "),
    @("E50", "  -1.98%  "),
    @("B51", "Frax"),
    @("C51", "https://coinranking.com/coin/KfWtaeV1W+frax-frax"),
    @("E51", "  -0.39%  ")
)

foreach ($u in $priceUpdates) {
    Set-TextValue $u[0] $u[1]
}

foreach ($u in $otherUpdates) {
    Set-PlainValue $u[0] $u[1]
}

